# Add Cart Valid and Invalid Features
# Rewrites the "Cart and Order Data Retrieval" sheet (3rd tab) so that the
# old Cart/Order test cases (TC_10 .. TC_15 content) are replaced with the
# new "Fetch All Carts" / "Fetch Single Cart with ID" / "Fetch Invalid
# Single Cart with ID" test cases (still labelled TC_10, TC_11, TC_12),
# and clears out the now-unused rows 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row: bold it ---
$ws.Range("A1:G1").Font.Bold = $true

# --- Row 2 : TC_10 - Fetch All Carts ---
$ws.Range("A2").Value = "TC_10"
$ws.Range("B2").Value = "Fetch All Carts"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "All carts appears in response"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "1. Send a GET request to /carts. <br> 2. Receive a JSON response with all carts items."

# --- Row 3 : TC_11 - Fetch Single Cart with ID ---
$ws.Range("A3").Value = "TC_11"
$ws.Range("B3").Value = "Fetch Single Cart with ID"
$ws.Range("C3").Value = "Valid Cart ID"
$ws.Range("D3").Value = "Single cart details appears in response"
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "1. Send a GET request to /carts/{cartID} with an valid cart ID. <br> 2.  Receive a JSON response with single Cart details"

# --- Row 4 : TC_12 - Fetch Invalid Single Cart with ID ---
$ws.Range("A4").Value = "TC_12"
$ws.Range("B4").Value = "Fetch Invalid Single Cart with ID"
$ws.Range("C4").Value = "Invalid Cart ID"
$ws.Range("D4").Value = "Error message and status code 404 "
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "1. Send a GET request to /carts/{cartID} with an invalid cart ID. <br> 2.  Receive an error message."

# --- Rows 5-7 : no longer used, wipe them but keep the blank formatted cells ---
$ws.Range("A5:D7").ClearContents()
$ws.Range("G5:G7").ClearContents()

# --- Highlight the Test Case ID column for the three rows (green fill) ---
$ws.Range("A2:A4").Interior.Color = 5287936

# --- Column sizing to match the refreshed layout ---
$ws.Columns.Item(2).ColumnWidth = 26.25
$ws.Columns.Item(5).Hidden = $true
$ws.Columns.Item(5).ColumnWidth = -0.8
$ws.Columns.Item(6).ColumnWidth = 21.92
$ws.Columns.Item(7).ColumnWidth = 102.92

# --- "Error Handling and Validation" sheet: widen the two description
#      columns (Excel "best fit" after adding the new rows on sheet 3
#      nudged the author to resize this sheet too) and move the selection ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Columns.Item(2).ColumnWidth = 18.25
$ws4.Columns.Item(3).ColumnWidth = 22.59
[void]$ws4.Range("B4").Select()

# --- Selection as left by the author after the edit (re-select sheet 3
#     last so it stays the active tab, matching the saved workbook view) ---
[void]$ws.Range("C13").Select()
